$d = $word.ActiveDocument

# --- Table indentation / cell-margin tweaks -------------------------------
$t = $d.Tables.Item(1)

# w:tblInd -5 dxa -> -10 dxa  (dxa / 20 = points)
$t.Rows.LeftIndent = -0.5

# w:tblCellMar left 103 dxa -> 98 dxa (table-level default cell margin)
$t.LeftPadding = 4.9

# w:tcMar left 103 dxa -> 98 dxa on every existing cell
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.LeftPadding = 4.9
    }
}

# --- Merge the two runs of the closing paragraph into a single run -------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "We combine the information from the recursive calls in linear time, thus from Eva Tardos, Algorithm Design, 5.10 the running time is O(nlogn).",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "We combine the information from the recursive calls in linear time, thus from Eva Tardos, Algorithm Design, 5.10 the running time is O(nlogn).",
    2) | Out-Null
